$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Concept Presentation paragraph: insert " to other team members," before
#    " during the actual presentation."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "slides during the actual presentation.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "slides to other team members, during the actual presentation.",
    2)

# ---------------------------------------------------------------------------
# 2) Implementation paragraph: add a trailing space after "Implementation:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Implementation:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Implementation: ",
    2)

# ---------------------------------------------------------------------------
# 3) Documentation paragraph: extend the closing sentence with the
#    additional work Vishal completed.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Implementation Evaluation sections of the report.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Implementation Evaluation sections. Furthermore I completed the Appendix 1 section. Finally I added my references and proof read the document, along with rest of the team, before we submitted it.",
    2)

# ---------------------------------------------------------------------------
# 4) Add a new "Overall:" paragraph right after the Documentation paragraph,
#    followed by a new blank paragraph.
# ---------------------------------------------------------------------------

# Locate the Documentation paragraph (the one that now contains the
# extended closing sentence we just wrote above).
$docPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Documentation:*") {
        $docPara = $para
    }
}

# Split a new, empty paragraph in right after the Documentation paragraph.
$insertionPoint = $d.Range($docPara.Range.End, $docPara.Range.End)
[void]$insertionPoint.InsertParagraphAfter()

# Re-locate that freshly created empty paragraph (it immediately follows the
# Documentation paragraph) and fill it in with the "Overall:" content.
$overallPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Start -eq $docPara.Range.End) {
        $overallPara = $para
    }
}

$overallXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:line="252" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">Overall: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:b/>
      <w:bCs/>
      <w:i/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>I feel like I was an active team member who contributed</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
'@
[void]$overallPara.Range.InsertXML($overallXml)

# ---------------------------------------------------------------------------
# 5) Add a new blank paragraph right after the "Overall:" paragraph.
# ---------------------------------------------------------------------------
$overallPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Overall:*") {
        $overallPara = $para
    }
}

$insertionPoint2 = $d.Range($overallPara.Range.End, $overallPara.Range.End)
[void]$insertionPoint2.InsertParagraphAfter()

$blankPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Start -eq $overallPara.Range.End) {
        $blankPara = $para
    }
}

$blankXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:line="252" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
'@
[void]$blankPara.Range.InsertXML($blankXml)
